# rdfization des gravures du MG
# - rename the "GMG" collection id to "MGE" (row 6, column A)
# - remove the now-obsolete "Partitions" row (row 7)
# - leave the selection on A6, matching the author's final cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: "GMG" -> "MGE"
$ws.Range("A6").Value = "MGE"

# Row 7 ("P" / "Partitions" / "Edition MEI" / "Partitions numérisées") is removed entirely
$ws.Rows(7).Delete()

# Match the saved selection state
$ws.Range("A6").Select() | Out-Null
